$d = $word.ActiveDocument

# Append: one blank paragraph, a heading-style paragraph, then a paragraph
# that ends up containing two separate runs of body text.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertAfter("My thoughts for the future (perhaps a paper)")

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertAfter("The socioeconomic models at the two different scale showed things like reversal of the directions of effects. Basically, in many cases there is no single relationships between the predictors and the response that is true across the whole country – the effects are really different for the different communes/provinces/regions. That’s why the models are shit. An idea – reanalyse the socioeconomic data, at the commune scale, BUT split the data by cluster. So either completely subset the data and do separate analyses for each cluster, OR remove the province random effect and instead use cluster. This will cluster the data into RE levels that actually match their socioeconomics. ")

# Start a new paragraph with the second chunk of text, then join it back
# onto the previous paragraph by deleting the intervening paragraph mark.
# This produces two independent <w:r> runs within a single <w:p>, matching
# how Word keeps runs separate when text was typed across an (undone) split.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertAfter("What I would expect to see are completely different effects (directions, sizes) between the different levels/clusters. When compared with the “global” effects, this would highlight the issue of scale.")

$paras = $d.Paragraphs
$joinIndex = $paras.Count - 1
$joinPara = $paras.Item($joinIndex)
$markStart = $joinPara.Range.End - 1
$markEnd = $joinPara.Range.End
$markRange = $d.Range($markStart, $markEnd)
$markRange.Delete()
